$wb = $excel.ActiveWorkbook

# --- Summary sheet: update value + selection ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F3").Value = 2986.88
$wsSummary.Activate()
$wsSummary.Range("B5").Select()

# --- Repayment Schedule sheet: update selection ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Activate()
$wsRepay.Range("C8").Select()

# --- Transactions sheet: update values + selection ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 217
$wsTrans.Range("A3").Value = 216
$wsTrans.Range("A4").Value = 214
$wsTrans.Range("A5").Value = 212
$wsTrans.Activate()
$wsTrans.Range("C4").Select()
